# Add a new "tntcoins3" sheet (a fresh snapshot copy of "tntcoins2") at the
# end of the workbook, and bump the "one4kat" player's stats to a newer
# snapshot.

$wb = $excel.ActiveWorkbook

# --- 1. Duplicate "tntcoins2" as a new trailing sheet named "tntcoins3" ---
$srcSheet = $wb.Worksheets.Item("tntcoins2")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$srcSheet.Copy($null, $lastSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "tntcoins3"

# --- 2. Update "one4kat" stats to the newer values ---
$ws = $wb.Worksheets.Item("one4kat")

# row 2 - available_layers
$ws.Range("B2").Value = 50
$ws.Range("C2").Value = 40
$ws.Range("E2").Value = 40
$ws.Range("G2").Value = 40
$ws.Range("I2").Value = 40

# row 3 - experience
$ws.Range("B3").Value = 588647.0833333324
$ws.Range("C3").Value = 15532.91666666663
$ws.Range("E3").Value = 15532.91666666663
$ws.Range("G3").Value = 15532.91666666663
$ws.Range("I3").Value = 15532.91666666663

# row 5 - coins
$ws.Range("B5").Value = 141312
$ws.Range("C5").Value = 2312
$ws.Range("E5").Value = 2312
$ws.Range("G5").Value = 2312
$ws.Range("I5").Value = 2312

# row 6 - damage_dealt
$ws.Range("B6").Value = 100822
$ws.Range("C6").Value = 2597
$ws.Range("E6").Value = 2597
$ws.Range("G6").Value = 2597
$ws.Range("I6").Value = 2597

# row 7 - deaths
$ws.Range("B7").Value = 3531
$ws.Range("C7").Value = 91
$ws.Range("E7").Value = 91
$ws.Range("G7").Value = 91
$ws.Range("I7").Value = 91

# row 9 - games_played
$ws.Range("B9").Value = 5934
$ws.Range("C9").Value = 160
$ws.Range("E9").Value = 160
$ws.Range("G9").Value = 160
$ws.Range("I9").Value = 160

# row 11 - sheep_thrown
$ws.Range("B11").Value = 56344
$ws.Range("C11").Value = 1097
$ws.Range("E11").Value = 1097
$ws.Range("G11").Value = 1097
$ws.Range("I11").Value = 1097

# row 13 - deaths_void
$ws.Range("B13").Value = 2586
$ws.Range("C13").Value = 67
$ws.Range("E13").Value = 67
$ws.Range("G13").Value = 67
$ws.Range("I13").Value = 67

# row 14 - wins
$ws.Range("B14").Value = 3906
$ws.Range("C14").Value = 135
$ws.Range("E14").Value = 135
$ws.Range("G14").Value = 135
$ws.Range("I14").Value = 135

# row 15 - kills
$ws.Range("B15").Value = 4211
$ws.Range("C15").Value = 112
$ws.Range("E15").Value = 112
$ws.Range("G15").Value = 112
$ws.Range("I15").Value = 112

# row 16 - kills_void
$ws.Range("B16").Value = 2763
$ws.Range("C16").Value = 66
$ws.Range("E16").Value = 66
$ws.Range("G16").Value = 66
$ws.Range("I16").Value = 66

# row 22 - playtime
$ws.Range("B22").Value = 622562
$ws.Range("C22").Value = 15403
$ws.Range("E22").Value = 15403
$ws.Range("G22").Value = 15403
$ws.Range("I22").Value = 15403
